$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$shp = $s.Shapes.Item(2)
$tf = $shp.TextFrame
$tr = $tf.TextRange

# bodyPr: normAutofit lnSpcReduction="10000" -> noAutofit
$tf.AutoSize = 0

# Paragraphs 1-6 keep their text, just shrink font size to 13pt (1300)
for ($i = 1; $i -le 6; $i++) {
    $tr.Paragraphs($i).Font.Size = 13
}

# Paragraph 7: "Функция search: ..." -> "Функция edit: ..."
$para7 = $tr.Paragraphs(7)
$para7.Text = ""
$para7.Text = "Функция edit: дублирует изменения в таблице в базу данных."
$para7.Font.Size = 13

# Paragraph 8 was empty (just endParaRPr) -> now holds the old "Функция search" text
$para8 = $tr.Paragraphs(8)
$para8.Text = "Функция search: производит поиск по названию товара в таблице."
$para8.Font.Size = 13
